# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" worksheet (fund holding detail) positioned
# between the existing "总计" and "2022-Q3" sheets, and adds a
# corresponding summary row to "总计", shifting the pre-existing
# 2022-Q3 / 2021-Q2 summary rows down by one.

$wb = $excel.ActiveWorkbook

# Helper: write a value that must be stored as TEXT even though it looks
# like a number (e.g. "012528" or "3.82"), without leaving a stray
# text-number-format style on the cell afterwards.
function Set-TextValue($cell, [string]$text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet — insert the new 2022-Q4 row at the top of
#    the data (row 2), pushing the existing 2022-Q3 / 2021-Q2 rows down.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A3:D3").Copy($total.Range("A4:D4"))
$total.Range("A2:D2").Copy($total.Range("A3:D3"))

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 2
$total.Cells.Item(2, 4).Value = 0.22

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2

# ---------------------------------------------------------------------
# 2) New "2022-Q4" detail sheet. Duplicate "2022-Q3" (via Copy) so the
#    new sheet inherits the same sheetPr / header styles / page margins,
#    then drop the extra template rows and overwrite with the real data.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# "2022-Q3" has 4 data rows (rows 2-5); "2022-Q4" only needs 2 (rows 2-3).
$q4.Rows.Item(5).Delete()
$q4.Rows.Item(4).Delete()

$q4.Cells.Item(2, 1).Value = 0
Set-TextValue $q4.Cells.Item(2, 2) "012528"
$q4.Cells.Item(2, 3).Value = "广发鑫睿一年持有期混合A"
Set-TextValue $q4.Cells.Item(2, 4) "3.82"
Set-TextValue $q4.Cells.Item(2, 5) "95.35"
Set-TextValue $q4.Cells.Item(2, 6) "3.84"
Set-TextValue $q4.Cells.Item(2, 7) "0.1467"
$q4.Cells.Item(2, 8).Value = 10

$q4.Cells.Item(3, 1).Value = 1
Set-TextValue $q4.Cells.Item(3, 2) "012529"
$q4.Cells.Item(3, 3).Value = "广发鑫睿一年持有期混合C"
Set-TextValue $q4.Cells.Item(3, 4) "1.87"
Set-TextValue $q4.Cells.Item(3, 5) "95.35"
Set-TextValue $q4.Cells.Item(3, 6) "3.84"
Set-TextValue $q4.Cells.Item(3, 7) "0.0718"
$q4.Cells.Item(3, 8).Value = 10

# Restore the original "last sheet is the active one" view state
# (the new sheet becomes active as a side effect of Copy/Add).
$wb.Worksheets.Item("2021-Q2").Activate()
